# Fill in the first six "CUT/CAD ..." sound-cue rows of the raccords
# table (column 2 was blank in the original upload) and drop the stale
# "_GoBack" last-edit-position bookmark that Word had left behind in an
# otherwise-empty cell further down the table.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$labels = @("CUT SON", "CAD MET SON", "CUT SON", "CAD FR MOV SON", "CAD SON", "CAD SON")

for ($i = 0; $i -lt $labels.Count; $i++) {
    $row = $t.Rows.Item($i + 1)
    $cell = $row.Cells.Item(2)
    $cell.Range.Text = $labels[$i]
}

# The "_GoBack" bookmark marks Word's last editing position; it sits
# alone in an otherwise-empty table cell. Find it, locate the cell it
# lives in, and clear that paragraph so the bookmark markup disappears
# along with it (mirrors Word's own cleanup of stale _GoBack marks).
$bookmark = $d.Bookmarks.Item("_GoBack")
$bookmarkCells = $bookmark.Range.Cells
$targetCell = $bookmarkCells.Item($bookmarkCells.Count)
$targetParagraph = $targetCell.Range.Paragraphs.Item(1)
$targetParagraph.Range.Delete()
